$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge the three runs "ratica " + "SUAP/SUE" + " n°" into a
# single run "ratica SUAP/SUE n°" (purely a run-merge, same visible text)
# ---------------------------------------------------------------------
$find1 = $d.Content.Find
$find1.ClearFormatting()
$found1 = $find1.Execute("ratica SUAP/SUE n" + [char]0x00B0, $false, $false, $false, $false, $false, $true, 1, $false, "ratica SUAP/SUE n" + [char]0x00B0, 2)
Write-Host "Step1 (merge ratica/SUAP/SUE/n run) found: $found1"

# ---------------------------------------------------------------------
# Change 2: remove the PEC Comando handoff -> the applicant no longer
# sends the request to the Comando's PEC, everything goes through
# SUAP/SUE/SUA instead.
#   a) trim the trailing ", inviando l'indirizzo di posta elettronica
#      certificata:" from the "Il richiedente..." paragraph, replacing
#      it with a final period.
#   b) delete the whole following paragraph that only contained the
#      <$PEC_COMANDO> merge field.
# ---------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.ClearFormatting()
$apostrophe = [char]0x2019
$target2 = ", inviando l" + $apostrophe + "indirizzo di posta elettronica certificata:"
$found2 = $find2.Execute($target2, $false, $false, $false, $false, $false, $true, 1, $false, ".", 2)
Write-Host "Step2a (trim PEC sentence) found: $found2"

$wdWithInTable = 12
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    $trimmed = $t.TrimEnd([char]13, [char]7)
    if ($trimmed -eq "<`$PEC_COMANDO>" -and -not $p.Range.Information($wdWithInTable)) {
        Write-Host "Step2b removing paragraph $i containing <`$PEC_COMANDO> (raw len $($t.Length))"
        $p.Range.Delete()
    }
}
